$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above row 11 (current "Description" row), shifting rows 11-21 down to 12-22.
$ws.Rows.Item(11).Insert()

# Fill in the new "Jurisdiction" row with an empty value cell.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the Date value (row 8, column B) to the new timestamp.
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
